# Add a "success" column (D) that flags, per distinct `list` value (rows
# 2-74, keyed by the running index in column A), whether that list is one
# of the first two observed lists (A = 1 or A = 2) -> "1", else "0".
# Matches commit "plot success rate 0-8".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: reuse the same bold/bordered/centered header look as B1/C1
# (style index 1 in the sheet) so it lands on the identical shared xf.
$ws.Range("D1").Value = "success"
$d1 = $ws.Range("D1")
$d1.Font.Bold = $true
$d1.HorizontalAlignment = -4108
$d1.VerticalAlignment = -4160
$d1.Borders.LineStyle = 1

# Data cells D2:D74 must hold TEXT "0"/"1" (not numbers) to match the
# source data (pandas-style boolean-as-string column). Format the range as
# text first so Excel stores the shared-string literal instead of
# auto-coercing "0"/"1" to numeric values, then restore the default
# (unstyled) look so no stray numeric cell format lingers on the cells.
$dataRange = $ws.Range("D2:D74")
$dataRange.NumberFormat = "@"

$rowCount = 74
for ($r = 2; $r -le $rowCount; $r++) {
    if ($r -eq 3 -or $r -eq 4) {
        $ws.Cells.Item($r, 4).Value = "1"
    } else {
        $ws.Cells.Item($r, 4).Value = "0"
    }
}

$dataRange.Style = "Normal"
